$d = $word.ActiveDocument
$d.Content.Find.Execute("54+1=", $true, $false, $false, $false, $false, $true, 1, $false, "94-5=", 2) | Out-Null
$d.Content.Find.Execute("56+33=", $true, $false, $false, $false, $false, $true, 1, $false, "66+3=", 2) | Out-Null
$d.Content.Find.Execute("78+7=", $true, $false, $false, $false, $false, $true, 1, $false, "2+87=", 2) | Out-Null
$d.Content.Find.Execute("34-1=", $true, $false, $false, $false, $false, $true, 1, $false, "89-65=", 2) | Out-Null
$d.Content.Find.Execute("35+15=", $true, $false, $false, $false, $false, $true, 1, $false, "83-60=", 2) | Out-Null
$d.Content.Find.Execute("85-61=", $true, $false, $false, $false, $false, $true, 1, $false, "18+46=", 2) | Out-Null
$d.Content.Find.Execute("27+43=", $true, $false, $false, $false, $false, $true, 1, $false, "33+8=", 2) | Out-Null
$d.Content.Find.Execute("1+60=", $true, $false, $false, $false, $false, $true, 1, $false, "66+13=", 2) | Out-Null
$d.Content.Find.Execute("30+35=", $true, $false, $false, $false, $false, $true, 1, $false, "9+46=", 2) | Out-Null
$d.Content.Find.Execute("46-28=", $true, $false, $false, $false, $false, $true, 1, $false, "22+40=", 2) | Out-Null
$d.Content.Find.Execute("70-4=", $true, $false, $false, $false, $false, $true, 1, $false, "3+81=", 2) | Out-Null
$d.Content.Find.Execute("68+4=", $true, $false, $false, $false, $false, $true, 1, $false, "24+33=", 2) | Out-Null
$d.Content.Find.Execute("79-3=", $true, $false, $false, $false, $false, $true, 1, $false, "89+4=", 2) | Out-Null
$d.Content.Find.Execute("81-30=", $true, $false, $false, $false, $false, $true, 1, $false, "26+37=", 2) | Out-Null
$d.Content.Find.Execute("59-10=", $true, $false, $false, $false, $false, $true, 1, $false, "66+26=", 2) | Out-Null
$d.Content.Find.Execute("50+6=", $true, $false, $false, $false, $false, $true, 1, $false, "6+76=", 2) | Out-Null
$d.Content.Find.Execute("77-30=", $true, $false, $false, $false, $false, $true, 1, $false, "46-11=", 2) | Out-Null
$d.Content.Find.Execute("66-15=", $true, $false, $false, $false, $false, $true, 1, $false, "7+18=", 2) | Out-Null
$d.Content.Find.Execute("28-3=", $true, $false, $false, $false, $false, $true, 1, $false, "53+35=", 2) | Out-Null
$d.Content.Find.Execute("59-30=", $true, $false, $false, $false, $false, $true, 1, $false, "3+43=", 2) | Out-Null
$d.Content.Find.Execute("16+2=", $true, $false, $false, $false, $false, $true, 1, $false, "70-42=", 2) | Out-Null
$d.Content.Find.Execute("19+42=", $true, $false, $false, $false, $false, $true, 1, $false, "77+13=", 2) | Out-Null
$d.Content.Find.Execute("59+29=", $true, $false, $false, $false, $false, $true, 1, $false, "33+42=", 2) | Out-Null
$d.Content.Find.Execute("86-21=", $true, $false, $false, $false, $false, $true, 1, $false, "14+35=", 2) | Out-Null
$d.Content.Find.Execute("43+6=", $true, $false, $false, $false, $false, $true, 1, $false, "75-6=", 2) | Out-Null
$d.Content.Find.Execute("77+21=", $true, $false, $false, $false, $false, $true, 1, $false, "82-45=", 2) | Out-Null
$d.Content.Find.Execute("62+2=", $true, $false, $false, $false, $false, $true, 1, $false, "3+90=", 2) | Out-Null
$d.Content.Find.Execute("84-38=", $true, $false, $false, $false, $false, $true, 1, $false, "48+9=", 2) | Out-Null
$d.Content.Find.Execute("45-5=", $true, $false, $false, $false, $false, $true, 1, $false, "94+5=", 2) | Out-Null
$d.Content.Find.Execute("63-14=", $true, $false, $false, $false, $false, $true, 1, $false, "97-14=", 2) | Out-Null
$d.Content.Find.Execute("58-26=", $true, $false, $false, $false, $false, $true, 1, $false, "88+4=", 2) | Out-Null
$d.Content.Find.Execute("86-28=", $true, $false, $false, $false, $false, $true, 1, $false, "65-2=", 2) | Out-Null
$d.Content.Find.Execute("94-3=", $true, $false, $false, $false, $false, $true, 1, $false, "78-20=", 2) | Out-Null
$d.Content.Find.Execute("10+36=", $true, $false, $false, $false, $false, $true, 1, $false, "18+2=", 2) | Out-Null
$d.Content.Find.Execute("6+64=", $true, $false, $false, $false, $false, $true, 1, $false, "46+26=", 2) | Out-Null
$d.Content.Find.Execute("87-61=", $true, $false, $false, $false, $false, $true, 1, $false, "23+15=", 2) | Out-Null
$d.Content.Find.Execute("69-26=", $true, $false, $false, $false, $false, $true, 1, $false, "31-27=", 2) | Out-Null
$d.Content.Find.Execute("54-7=", $true, $false, $false, $false, $false, $true, 1, $false, "51-27=", 2) | Out-Null
$d.Content.Find.Execute("82-46=", $true, $false, $false, $false, $false, $true, 1, $false, "85-45=", 2) | Out-Null
$d.Content.Find.Execute("63+6=", $true, $false, $false, $false, $false, $true, 1, $false, "8+70=", 2) | Out-Null
$d.Content.Find.Execute("10+60=", $true, $false, $false, $false, $false, $true, 1, $false, "15+58=", 2) | Out-Null
$d.Content.Find.Execute("45-2=", $true, $false, $false, $false, $false, $true, 1, $false, "96-6=", 2) | Out-Null
$d.Content.Find.Execute("25+37=", $true, $false, $false, $false, $false, $true, 1, $false, "67-62=", 2) | Out-Null
$d.Content.Find.Execute("60-15=", $true, $false, $false, $false, $false, $true, 1, $false, "27+23=", 2) | Out-Null
$d.Content.Find.Execute("77-1=", $true, $false, $false, $false, $false, $true, 1, $false, "49+41=", 2) | Out-Null
$d.Content.Find.Execute("56-5=", $true, $false, $false, $false, $false, $true, 1, $false, "54+42=", 2) | Out-Null
$d.Content.Find.Execute("48+37=", $true, $false, $false, $false, $false, $true, 1, $false, "29+70=", 2) | Out-Null
$d.Content.Find.Execute("87-36=", $true, $false, $false, $false, $false, $true, 1, $false, "70-28=", 2) | Out-Null
$d.Content.Find.Execute("62-2=", $true, $false, $false, $false, $false, $true, 1, $false, "68+8=", 2) | Out-Null
$d.Content.Find.Execute("44-41=", $true, $false, $false, $false, $false, $true, 1, $false, "42+18=", 2) | Out-Null
$d.Content.Find.Execute("16+38=", $true, $false, $false, $false, $false, $true, 1, $false, "89-59=", 2) | Out-Null
$d.Content.Find.Execute("25+68=", $true, $false, $false, $false, $false, $true, 1, $false, "9+50=", 2) | Out-Null
$d.Content.Find.Execute("42-17=", $true, $false, $false, $false, $false, $true, 1, $false, "34+5=", 2) | Out-Null
$d.Content.Find.Execute("21-3=", $true, $false, $false, $false, $false, $true, 1, $false, "85-82=", 2) | Out-Null
$d.Content.Find.Execute("47+25=", $true, $false, $false, $false, $false, $true, 1, $false, "66+12=", 2) | Out-Null
$d.Content.Find.Execute("31-11=", $true, $false, $false, $false, $false, $true, 1, $false, "78-61=", 2) | Out-Null
$d.Content.Find.Execute("33+52=", $true, $false, $false, $false, $false, $true, 1, $false, "67+13=", 2) | Out-Null
$d.Content.Find.Execute("42+29=", $true, $false, $false, $false, $false, $true, 1, $false, "37+36=", 2) | Out-Null
$d.Content.Find.Execute("66+23=", $true, $false, $false, $false, $false, $true, 1, $false, "57-45=", 2) | Out-Null
$d.Content.Find.Execute("43-2=", $true, $false, $false, $false, $false, $true, 1, $false, "81+8=", 2) | Out-Null
$d.Content.Find.Execute("1+37=", $true, $false, $false, $false, $false, $true, 1, $false, "26+37=", 2) | Out-Null
$d.Content.Find.Execute("0+41=", $true, $false, $false, $false, $false, $true, 1, $false, "11+23=", 2) | Out-Null
$d.Content.Find.Execute("83-45=", $true, $false, $false, $false, $false, $true, 1, $false, "46+38=", 2) | Out-Null
$d.Content.Find.Execute("51+15=", $true, $false, $false, $false, $false, $true, 1, $false, "9+66=", 2) | Out-Null
$d.Content.Find.Execute("99-6=", $true, $false, $false, $false, $false, $true, 1, $false, "0+83=", 2) | Out-Null
$d.Content.Find.Execute("86-45=", $true, $false, $false, $false, $false, $true, 1, $false, "89-69=", 2) | Out-Null
$d.Content.Find.Execute("37-10=", $true, $false, $false, $false, $false, $true, 1, $false, "78-11=", 2) | Out-Null
$d.Content.Find.Execute("27+62=", $true, $false, $false, $false, $false, $true, 1, $false, "25+64=", 2) | Out-Null
$d.Content.Find.Execute("64-57=", $true, $false, $false, $false, $false, $true, 1, $false, "52+40=", 2) | Out-Null
$d.Content.Find.Execute("56-14=", $true, $false, $false, $false, $false, $true, 1, $false, "59+13=", 2) | Out-Null
$d.Content.Find.Execute("4+11=", $true, $false, $false, $false, $false, $true, 1, $false, "48-3=", 2) | Out-Null
$d.Content.Find.Execute("27-17=", $true, $false, $false, $false, $false, $true, 1, $false, "33+3=", 2) | Out-Null
$d.Content.Find.Execute("74-27=", $true, $false, $false, $false, $false, $true, 1, $false, "16+39=", 2) | Out-Null
$d.Content.Find.Execute("71-7=", $true, $false, $false, $false, $false, $true, 1, $false, "38+8=", 2) | Out-Null
$d.Content.Find.Execute("98-22=", $true, $false, $false, $false, $false, $true, 1, $false, "76-44=", 2) | Out-Null
$d.Content.Find.Execute("56+37=", $true, $false, $false, $false, $false, $true, 1, $false, "91-42=", 2) | Out-Null
$d.Content.Find.Execute("66+6=", $true, $false, $false, $false, $false, $true, 1, $false, "72-56=", 2) | Out-Null
$d.Content.Find.Execute("46-38=", $true, $false, $false, $false, $false, $true, 1, $false, "69+9=", 2) | Out-Null
$d.Content.Find.Execute("32+53=", $true, $false, $false, $false, $false, $true, 1, $false, "5+68=", 2) | Out-Null
$d.Content.Find.Execute("22-16=", $true, $false, $false, $false, $false, $true, 1, $false, "78-60=", 2) | Out-Null
$d.Content.Find.Execute("67-8=", $true, $false, $false, $false, $false, $true, 1, $false, "29+23=", 2) | Out-Null
$d.Content.Find.Execute("29+48=", $true, $false, $false, $false, $false, $true, 1, $false, "9+26=", 2) | Out-Null
$d.Content.Find.Execute("83+7=", $true, $false, $false, $false, $false, $true, 1, $false, "59-17=", 2) | Out-Null
$d.Content.Find.Execute("18+0=", $true, $false, $false, $false, $false, $true, 1, $false, "60-43=", 2) | Out-Null
$d.Content.Find.Execute("70+23=", $true, $false, $false, $false, $false, $true, 1, $false, "46+43=", 2) | Out-Null
$d.Content.Find.Execute("76+18=", $true, $false, $false, $false, $false, $true, 1, $false, "21+60=", 2) | Out-Null
$d.Content.Find.Execute("17+71=", $true, $false, $false, $false, $false, $true, 1, $false, "72+11=", 2) | Out-Null
$d.Content.Find.Execute("41+26=", $true, $false, $false, $false, $false, $true, 1, $false, "54-44=", 2) | Out-Null
$d.Content.Find.Execute("63-45=", $true, $false, $false, $false, $false, $true, 1, $false, "33-14=", 2) | Out-Null
$d.Content.Find.Execute("25-21=", $true, $false, $false, $false, $false, $true, 1, $false, "61+21=", 2) | Out-Null
$d.Content.Find.Execute("79-26=", $true, $false, $false, $false, $false, $true, 1, $false, "77-69=", 2) | Out-Null
$d.Content.Find.Execute("48-33=", $true, $false, $false, $false, $false, $true, 1, $false, "18-17=", 2) | Out-Null
$d.Content.Find.Execute("99-97=", $true, $false, $false, $false, $false, $true, 1, $false, "95-86=", 2) | Out-Null
$d.Content.Find.Execute("71+21=", $true, $false, $false, $false, $false, $true, 1, $false, "31-23=", 2) | Out-Null
$d.Content.Find.Execute("3+74=", $true, $false, $false, $false, $false, $true, 1, $false, "79-44=", 2) | Out-Null
$d.Content.Find.Execute("81-66=", $true, $false, $false, $false, $false, $true, 1, $false, "67+9=", 2) | Out-Null
$d.Content.Find.Execute("48-5=", $true, $false, $false, $false, $false, $true, 1, $false, "60+7=", 2) | Out-Null
$d.Content.Find.Execute("30+51=", $true, $false, $false, $false, $false, $true, 1, $false, "98-71=", 2) | Out-Null
$d.Content.Find.Execute("36+63=", $true, $false, $false, $false, $false, $true, 1, $false, "96-88=", 2) | Out-Null
$d.Content.Find.Execute("44-32=", $true, $false, $false, $false, $false, $true, 1, $false, "56-50=", 2) | Out-Null
